# Daily attendance processing - 2025-12-25 07:35:19
# Normalize the "Recorded By" (column G) author-list ordering: swap the
# first two comma-separated names/emails in each cell (trailing entries,
# if any, are left in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("G1:G157")

$rng.Replace("system, backup@backdoor.com, System", "backup@backdoor.com, system, System")
$rng.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")
$rng.Replace("System, admin@admin.com", "admin@admin.com, System")
$rng.Replace("dnasr281@gmail.com, admin@admin.com", "admin@admin.com, dnasr281@gmail.com")
